$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCategories = $wb.Worksheets.Item("Categories")

# Update the row_id label text on the "Variables" sheet so it now reads
# "Unique identifier for the row in Opal" instead of the old
# "Unique identifier for the row". Excel will drop the now-unused shared
# string and append the new one at the end of the shared string table,
# which re-indexes every other <v> string reference that follows it.
$wsVariables.Range("D2").Value = "Unique identifier for the row in Opal"

# "Categories" keeps a lingering selection of C8 (no longer the active
# tab, but its sheetView remembers the last selection made there) - set
# this first so it doesn't clobber the final active-tab state below.
$wsCategories.Activate()
$wsCategories.Range("C8").Select()

# Switch the active sheet/tab + selection: "Variables" becomes the
# selected tab (was "Categories"), with D2 as the active cell.
$wsVariables.Activate()
$wsVariables.Range("D2").Select()
